$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Re-sort the account-statement rows (16-19) so they are grouped by
# "Periodo Mora" (1702 first, then 1703) instead of by worker.
$ws.Range("E16").Value = "1702"
$ws.Range("C17").Value = "1007210721"
$ws.Range("D17").Value = "JORGE LUIS ROENES CASTILLO"
$ws.Range("E17").Value = "1702"
$ws.Range("F17").Value = 30000
$ws.Range("G17").Value = 750000
$ws.Range("C18").Value = "1235039833"
$ws.Range("D18").Value = "LISBETH PAOLA CUADRO MORALES"
$ws.Range("E18").Value = "1703"
$ws.Range("F18").Value = 29509
$ws.Range("G18").Value = 737717
$ws.Range("E19").Value = "1703"
